$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 with the new group name and new course description,
# mirroring the content and style used in the existing rows.
$ws.Range("A3").Value = "Grupo de Investigación en Recursos Hídricos y Saneamiento Ambiental"
$ws.Range("B3").Value = "5.- Nombre del Curso: Monitoreo y Evaluación de la Calidad del Agua Fecha acto administrativo curso: 2017-12-05 Número acto administrativo curso: 318 Programa académico: Maestría en Ingeniería Civil"

# Copy the formatting from A2 (bold font, thin border, centered
# horizontally, top vertical alignment - same style used for A1/A2) onto A3.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
